$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Data")

# 1. Move the existing "Poudre Wilderness Volunteers" row (row 2) down to a
#    new row 7, preserving its formatting (styles + row height).
$ws.Range("A2:L2").Copy($ws.Range("A7:L7"))
$ws.Range("K7").ClearContents()
$ws.Rows.Item(7).RowHeight = 15

# 2. Overwrite row 2 with the new "Cameron Pass Nordic Ranger Program" entry
#    (plain / default formatting, like a freshly typed row).
$ws.Range("A2:L2").Clear()
$ws.Range("A2").Value = "Cameron Pass Nordic Ranger Program"
$ws.Range("B2").Value = "Volunteer"
$ws.Range("C2").Value = "Nordic trail maintenance"
$ws.Range("D2").Value = "Nordic trail maintenance, education."
$ws.Range("E2").Value = "https://www.fs.usda.gov/detail/arp/workingtogether/volunteering/?cid=stelprdb5213201"
$ws.Range("F2").Value = "Yes"
$ws.Range("G2").Value = -105.882266
$ws.Range("H2").Value = 40.527357000000002
$ws.Range("L2").Value = "Coordinates for Moose Visitor Center"
$ws.Rows.Item(2).AutoFit()

# 3. Column width tweaks for columns A and E.
$ws.Columns.Item(1).ColumnWidth = 31.666666666666668
$ws.Columns.Item(5).ColumnWidth = 75.83333333333333

# 4. Update the active selection.
$ws.Range("H3").Select() | Out-Null
